$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text (not auto-converted to a number),
# then restore its original (default/"Normal") cell style so no stray
# cell-level number formatting is introduced.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Enterprises density (per 1000 people) - row 13
Set-TextValue "B13" "37.76"
Set-TextValue "C13" "3.18"
Set-TextValue "D13" "40.94"

# Employment (% of total) - row 14
Set-TextValue "B14" "37.07"
Set-TextValue "C14" "44.77"
Set-TextValue "D14" "81.84"

# Enterprises (% of total) - row 16
Set-TextValue "B16" "92.07"
Set-TextValue "C16" "7.75"
Set-TextValue "D16" "99.82"

# Value added to the economy (% of total) - row 20
Set-TextValue "B20" "26.57"
Set-TextValue "C20" "49.71"
Set-TextValue "D20" "76.28"
